$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(2)

# Locate the " configured logging" text that ends the "SLF4J Logger"
# bullet so we don't depend on a hard-coded character offset.
$tr = $sh.TextFrame.TextRange
$full = $tr.Text
$needle = " configured logging"
$idx = $full.IndexOf($needle)

# Shrink that run down to " configured " (drop the trailing "logging");
# it gets re-added as its own run immediately below, splitting the
# original run into two runs just like the authored edit.
$target = $tr.Characters($idx + 1, $needle.Length)
$target.Text = " configured "

# Re-insert "logging" as a new trailing run in the same paragraph.
$tr2 = $sh.TextFrame.TextRange
$tr2.InsertAfter("logging") | Out-Null

# Add the new "Junit: For unit testing" paragraph, split across two
# runs just like the authored edit.
$tr3 = $sh.TextFrame.TextRange
$tr3.InsertAfter("`rJunit: For ") | Out-Null

$tr4 = $sh.TextFrame.TextRange
$tr4.InsertAfter("unit testing") | Out-Null
